# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Updates the metadata sheet so that:
#  - fondo-de-contingencia, estado-de-la-informacion and tipo-de-presupuesto
#    become "iaest-measure:" (type "medida" / "xsd:int") instead of
#    "iaest-dimension:" (type "dim" / "skos:Concept") columns, and lose
#    their mapping-*.xlsx files.
#  - municipio-nombre becomes a "sdmx-dimension:refArea" dimension (type
#    "dim" / "URI-Municipio") instead of an "iaest-measure:" measure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - concept URIs
$ws.Range("C2").Value = "iaest-measure:fondo-de-contingencia"
$ws.Range("E2").Value = "iaest-measure:estado-de-la-informacion"
$ws.Range("M2").Value = "sdmx-dimension:refArea"
$ws.Range("O2").Value = "iaest-measure:tipo-de-presupuesto"

# Row 3 - "medida" / "dim"
$ws.Range("C3").Value = "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("M3").Value = "dim"
$ws.Range("O3").Value = "medida"

# Row 4 - "xsd:int" / "skos:Concept" / "URI-*"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("M4").Value = "URI-Municipio"
$ws.Range("O4").Value = "xsd:int"

# Row 5 - mapping files no longer needed for the curated columns
$ws.Range("C5").Clear()
$ws.Range("E5").Clear()
$ws.Range("O5").Clear()
